$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor name (short form) in the header block
$ws.Range("G6").Value = "P.S"

# Supervisor sign-off: name and date
$ws.Range("A27").Value = "Prakruti Sinha"
$ws.Range("D27").Value = Get-Date -Year 2014 -Month 2 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("D27").NumberFormat = $ws.Range("D25").NumberFormat
